$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(14, 15, 334.7633794933903, 0.0416667, 0.0416667, 0.0833334, 502.1450692400854, 0.3, 0.4, 0.11, 0.37),
    @(15, 16, 286.7480913238156, 0.0138889, 0.0138889, 0.0277778, 430.1221369857234, 0.33, 0.49, 0.12, 0.36),
    @(16, 17, 262.4045378156197, 0, 0, 0, 393.6068067234295, 0.39, 0.46, 0.14, 0.39),
    @(17, 18, 309.1354664251031, 0, 0, 0, 463.7031996376547, 0.37, 0.48, 0.19, 0.34),
    @(18, 19, 384.636049180636, 0, 0, 0, 576.954073770954, 0.35, 0.4, 0.18, 0.39),
    @(19, 20, 433.5207499755035, 0, 0, 0, 650.2811249632553, 0.38, 0.41, 0.17, 0.33),
    @(20, 21, 393.8044005019047, 0, 0, 0, 590.7066007528571, 0.35, 0.49, 0.16, 0.37),
    @(21, 22, 333.3011855326704, 0, 0, 0, 499.9517782990056, 0.3, 0.4, 0.19, 0.34),
    @(22, 23, 282.2034344188766, 0, 0, 0, 423.3051516283149, 0.39, 0.41, 0.14, 0.3),
    @(23, 24, 208.1255268683671, 0, 0, 0, 312.1882903025506, 0.32, 0.45, 0.18, 0.39),
    @(24, 25, 165.1488800499209, 0, 0, 0, 247.7233200748814, 0.4, 0.48, 0.2, 0.4),
    @(25, 26, 130.2933375268219, 0, 0, 0, 195.4400062902328, 0.37, 0.46, 0.1, 0.33),
    @(26, 27, 103.6774555661557, 0, 0, 0, 155.5161833492336, 0.38, 0.48, 0.19, 0.32),
    @(27, 28, 97.78916096758275, 0, 0, 0, 146.6837414513741, 0.35, 0.43, 0.14, 0.32),
    @(28, 29, 94.78573553475324, 0, 0, 0, 142.1786033021299, 0.4, 0.49, 0.14, 0.39),
    @(29, 30, 98.71785172641819, 0, 0, 0, 148.0767775896273, 0.37, 0.46, 0.15, 0.33)
)

# Reference cell that already carries the bold/centered/bordered style used
# throughout column A of the table (A2 in the pre-edit sheet).
$styleRef = $ws.Cells.Item(2, 1)

for ($i = 0; $i -lt $data.Length; $i++) {
    $r = $i + 2
    $row = $data[$i]
    for ($j = 0; $j -lt $row.Length; $j++) {
        $c = $j + 1
        $ws.Cells.Item($r, $c).Value = $row[$j]
    }
}

# Re-apply the column-A style (bold/centered/bordered) to every data row,
# including the five newly-appended rows, by copying formats from the
# reference cell.
$styleRef.Copy() | Out-Null
$ws.Range("A2:A17").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0
